# Swap the B and E:AD column contents between paired rows.
# (Column A keeps the running index, and C/D are identical within each
# pair so they do not need to be touched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(25, 26),
    @(50, 51),
    @(99, 100),
    @(103, 104),
    @(129, 130),
    @(143, 144),
    @(148, 149)
)

$columns = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $columns) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
